# Test different types of emotions
# Update the "EMOTION" and "EVENT" values, and their related numeric
# readings, for the single data row (row 2) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# MOOD value (A2)
$ws.Range("A2").Value = 0.8896626830101013

# EMOTION (B2): Shame -> Gloating
$ws.Range("B2").Value = "Gloating"

# INTENSITY value (C2)
$ws.Range("C2").Value = 2.8645248413085938

# EVENT (D2): Fly -> BecomeRich
$ws.Range("D2").Value = "BecomeRich"
